$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold textual data (prices / percentages) that must remain
# plain text even though many values look numeric (e.g. "569.95", "0.355").
# Temporarily force the affected range to Text format so Excel does not
# auto-convert the assigned strings into numbers, then restore the original
# (default/"Normal") style once all values have been written.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '63.726.80'
$ws.Range('E2').Value = '  +5.04%  '
$ws.Range('D3').Value = '2.482.47'
$ws.Range('E3').Value = '  +6.11%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '569.95'
$ws.Range('E5').Value = '  +4.05%  '
$ws.Range('D6').Value = '143.74'
$ws.Range('E6').Value = '  +9.42%  '
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('E8').Value = '  +2.49%  '
$ws.Range('D9').Value = '2.481.34'
$ws.Range('E9').Value = '  +6.31%  '
$ws.Range('E10').Value = '  +4.41%  '
$ws.Range('D11').Value = '5.75'
$ws.Range('E11').Value = '  +4.44%  '
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = '0.355'
$ws.Range('D14').Value = '26.53'
$ws.Range('E14').Value = '  +12.62%  '
$ws.Range('D15').Value = '2.926.20'
$ws.Range('E15').Value = '  +5.78%  '
$ws.Range('D16').Value = '63.558.45'
$ws.Range('E16').Value = '  +4.83%  '
$ws.Range('E17').Value = '  +6.75%  '
$ws.Range('D18').Value = '2.483.44'
$ws.Range('E18').Value = '  +5.93%  '
$ws.Range('D19').Value = '11.35'
$ws.Range('E19').Value = '  +6.78%  '
$ws.Range('D20').Value = '342.98'
$ws.Range('E20').Value = '  +8.90%  '
$ws.Range('D21').Value = '4.33'
$ws.Range('E21').Value = '  +5.95%  '
$ws.Range('E22').Value = '  +4.18%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '65.93'
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  +9.29%  '
$ws.Range('D28').Value = '8.26'
$ws.Range('E28').Value = '  +4.37%  '
$ws.Range('E29').Value = '  +6.79%  '
$ws.Range('D30').Value = '0.0₃0830'
$ws.Range('E30').Value = '  +13.13%  '
$ws.Range('E31').Value = '  +15.42%  '
$ws.Range('E32').Value = '  +8.26%  '
$ws.Range('D33').Value = '177.37'
$ws.Range('E33').Value = '  +3.01%  '
$ws.Range('E34').Value = '  +11.10%  '
$ws.Range('E35').Value = '  +4.26%  '
$ws.Range('D36').Value = '19.04'
$ws.Range('E36').Value = '  +5.85%  '
$ws.Range('D37').Value = '372.86'
$ws.Range('E37').Value = '  +14.62%  '
$ws.Range('D38').Value = '4.49'
$ws.Range('E38').Value = '  +8.24%  '
$ws.Range('D40').Value = '1.72'
$ws.Range('E40').Value = '  +12.35%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = '40.41'
$ws.Range('E42').Value = '  +5.51%  '
$ws.Range('D43').Value = '151.57'
$ws.Range('E43').Value = '  +10.66%  '
$ws.Range('E44').Value = '  +6.85%  '
$ws.Range('D45').Value = '20.92'
$ws.Range('E45').Value = '  +9.15%  '
$ws.Range('D46').Value = '0.602'
$ws.Range('E46').Value = '  +5.11%  '
$ws.Range('D47').Value = '0.0969'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('E48').Value = '  +6.00%  '
$ws.Range('E49').Value = '  +8.36%  '
$ws.Range('E50').Value = '  +5.20%  '
$ws.Range('D51').Value = '18.20'
$ws.Range('E51').Value = '  +6.94%  '

$dataRange.Style = "Normal"
